$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 97.5  # H11
$ws.Cells.Item(11, 9).Value = 97.5  # I11
$ws.Cells.Item(11, 11).Value = 97.5  # K11
$ws.Cells.Item(11, 13).Value = 42.5  # M11
# Row 32
$ws.Cells.Item(32, 8).Value = 2325.5  # H32
$ws.Cells.Item(32, 9).Value = 0  # I32
$ws.Cells.Item(32, 11).Value = 0  # K32
$ws.Cells.Item(32, 13).ClearContents()  # M32 was -674
# Row 129
$ws.Cells.Item(129, 8).Value = 1223.7354  # H129
$ws.Cells.Item(129, 9).Value = 512.7143  # I129
$ws.Cells.Item(129, 10).Value = 1408.0741  # J129
$ws.Cells.Item(129, 11).Value = 1538.1429  # K129
$ws.Cells.Item(129, 12).Value = 4224.2223  # L129
$ws.Cells.Item(129, 13).Value = 3461.8571  # M129
$ws.Cells.Item(129, 14).Value = -14224.2223  # N129
# Row 137
$ws.Cells.Item(137, 8).Value = 1714  # H137
$ws.Cells.Item(137, 9).Value = 1823.0344  # I137
$ws.Cells.Item(137, 10).Value = 923.5  # J137
$ws.Cells.Item(137, 11).Value = 5469.1032  # K137
$ws.Cells.Item(137, 12).Value = 2770.5  # L137
$ws.Cells.Item(137, 13).Value = -2919.1032  # M137
$ws.Cells.Item(137, 14).Value = -7870.5  # N137

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 529861.8  # H32
$ws.Cells.Item(32, 9).Value = 568853.25  # I32
$ws.Cells.Item(32, 10).Value = 29471.334  # J32
$ws.Cells.Item(32, 11).Value = 568853.25  # K32
$ws.Cells.Item(32, 12).Value = 29471.334  # L32
$ws.Cells.Item(32, 13).Value = -568566.25  # M32
$ws.Cells.Item(32, 14).Value = -30045.334  # N32
# Row 110
$ws.Cells.Item(110, 8).Value = 115013.875  # H110
$ws.Cells.Item(110, 9).Value = 152351.83  # I110
$ws.Cells.Item(110, 11).Value = 152351.83  # K110
$ws.Cells.Item(110, 13).Value = -150306.83  # M110
# Row 123
$ws.Cells.Item(123, 8).Value = 24428  # H123
$ws.Cells.Item(123, 10).Value = 24428  # J123
$ws.Cells.Item(123, 12).Value = 24428  # L123
$ws.Cells.Item(123, 14).Value = -34228  # N123
# Row 132
$ws.Cells.Item(132, 8).Value = 2136.606  # H132
$ws.Cells.Item(132, 9).Value = 1241.8909  # I132
$ws.Cells.Item(132, 11).Value = 3725.6727  # K132
$ws.Cells.Item(132, 13).Value = -1195.6727  # M132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Cells.Item(8, 8).Value = 3126  # H8
$ws.Cells.Item(8, 9).Value = 3901.3333  # I8
$ws.Cells.Item(8, 10).Value = 800  # J8
$ws.Cells.Item(8, 11).Value = 3901.3333  # K8
$ws.Cells.Item(8, 12).Value = 800  # L8
$ws.Cells.Item(8, 13).Value = -3761.3333  # M8
$ws.Cells.Item(8, 14).Value = -1080  # N8 new cell
# Row 22
$ws.Cells.Item(22, 8).Value = 291  # H22
$ws.Cells.Item(22, 9).Value = 291  # I22
$ws.Cells.Item(22, 11).Value = 291  # K22
$ws.Cells.Item(22, 13).Value = -118  # M22
# Row 80
$ws.Cells.Item(80, 8).Value = 1209.6666  # H80
$ws.Cells.Item(80, 10).Value = 203.1  # J80
$ws.Cells.Item(80, 12).Value = 203.1  # L80
$ws.Cells.Item(80, 14).Value = -2199.1  # N80
# Row 83
$ws.Cells.Item(83, 8).Value = 1209.6666  # H83
$ws.Cells.Item(83, 10).Value = 203.1  # J83
$ws.Cells.Item(83, 12).Value = 1015.5  # L83
$ws.Cells.Item(83, 14).Value = -10999.5  # N83
# Row 105
$ws.Cells.Item(105, 8).Value = 7815202  # H105
$ws.Cells.Item(105, 9).Value = 8336082  # I105
$ws.Cells.Item(105, 11).Value = 8336082  # K105
$ws.Cells.Item(105, 13).Value = -8334335  # M105

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 4288.217  # H31
$ws.Cells.Item(31, 9).Value = 1105.7778  # I31
$ws.Cells.Item(31, 10).Value = 6892.0303  # J31
$ws.Cells.Item(31, 11).Value = 1105.7778  # K31
$ws.Cells.Item(31, 12).Value = 6892.0303  # L31
$ws.Cells.Item(31, 13).Value = -810.7778000000001  # M31
$ws.Cells.Item(31, 14).Value = -7482.0303  # N31
# Row 34
$ws.Cells.Item(34, 8).Value = 4288.217  # H34
$ws.Cells.Item(34, 9).Value = 1105.7778  # I34
$ws.Cells.Item(34, 10).Value = 6892.0303  # J34
$ws.Cells.Item(34, 11).Value = 1105.7778  # K34
$ws.Cells.Item(34, 12).Value = 6892.0303  # L34
$ws.Cells.Item(34, 13).Value = -903.7778000000001  # M34
$ws.Cells.Item(34, 14).Value = -7296.0303  # N34
# Row 39
$ws.Cells.Item(39, 8).Value = 0  # H39
$ws.Cells.Item(39, 9).Value = 0  # I39
$ws.Cells.Item(39, 11).Value = 0  # K39
$ws.Cells.Item(39, 13).ClearContents()  # M39 was -4609
# Row 49
$ws.Cells.Item(49, 8).Value = 0  # H49
$ws.Cells.Item(49, 9).Value = 0  # I49
$ws.Cells.Item(49, 11).Value = 0  # K49
$ws.Cells.Item(49, 13).ClearContents()  # M49 was -4818
# Row 132
$ws.Cells.Item(132, 8).Value = 43211650  # H132
$ws.Cells.Item(132, 10).Value = 18520506  # J132
$ws.Cells.Item(132, 12).Value = 55561518  # L132
$ws.Cells.Item(132, 14).Value = -55566578  # N132

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Cells.Item(109, 8).Value = 3660.842  # H109
$ws.Cells.Item(109, 9).Value = 847.3333  # I109
$ws.Cells.Item(109, 11).Value = 2541.9999  # K109
$ws.Cells.Item(109, 13).Value = -1501.9999  # M109
# Row 120
$ws.Cells.Item(120, 8).Value = 10942.941  # H120
$ws.Cells.Item(120, 9).Value = 7257.5  # I120
$ws.Cells.Item(120, 10).Value = 12076.923  # J120
$ws.Cells.Item(120, 11).Value = 21772.5  # K120
$ws.Cells.Item(120, 12).Value = 36230.769  # L120
$ws.Cells.Item(120, 13).Value = -16934.5  # M120
$ws.Cells.Item(120, 14).Value = -45906.769  # N120
# Row 131
$ws.Cells.Item(131, 8).Value = 4777.5806  # H131
$ws.Cells.Item(131, 9).Value = 457.22223  # I131
$ws.Cells.Item(131, 10).Value = 6545  # J131
$ws.Cells.Item(131, 11).Value = 1371.66669  # K131
$ws.Cells.Item(131, 12).Value = 19635  # L131
$ws.Cells.Item(131, 13).Value = 3668.33331  # M131
$ws.Cells.Item(131, 14).Value = -29715  # N131

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Cells.Item(32, 8).Value = 35000  # H32
$ws.Cells.Item(32, 10).Value = 35000  # J32
$ws.Cells.Item(32, 12).Value = 35000  # L32
$ws.Cells.Item(32, 14).Value = -35592  # N32
# Row 107
$ws.Cells.Item(107, 8).Value = 0  # H107
$ws.Cells.Item(107, 9).Value = 0  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 0  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).ClearContents()  # M107 was 1644.8
$ws.Cells.Item(107, 14).ClearContents()  # N107 was -4138
# Row 113
$ws.Cells.Item(113, 8).Value = 43750.582  # H113
$ws.Cells.Item(113, 9).Value = 51968.7  # I113
$ws.Cells.Item(113, 11).Value = 51968.7  # K113
$ws.Cells.Item(113, 13).Value = -49798.7  # M113
# Row 132
$ws.Cells.Item(132, 8).Value = 2904.3455  # H132
$ws.Cells.Item(132, 9).Value = 2532.907  # I132
$ws.Cells.Item(132, 11).Value = 7598.721  # K132
$ws.Cells.Item(132, 13).Value = -5068.721  # M132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Cells.Item(24, 8).Value = 70007  # H24
$ws.Cells.Item(24, 9).Value = 0  # I24
$ws.Cells.Item(24, 11).Value = 0  # K24
$ws.Cells.Item(24, 13).ClearContents()  # M24 was -663
# Row 116
$ws.Cells.Item(116, 8).Value = 39800  # H116
$ws.Cells.Item(116, 10).Value = 39800  # J116
$ws.Cells.Item(116, 12).Value = 39800  # L116
$ws.Cells.Item(116, 14).Value = -48978  # N116 new cell
# Row 132
$ws.Cells.Item(132, 8).Value = 2694.244  # H132
$ws.Cells.Item(132, 9).Value = 2387.7856  # I132
$ws.Cells.Item(132, 11).Value = 7163.3568  # K132
$ws.Cells.Item(132, 13).Value = -4633.3568  # M132
# Row 136
$ws.Cells.Item(136, 8).Value = 8773843  # H136
$ws.Cells.Item(136, 9).Value = 1932.1666  # I136
$ws.Cells.Item(136, 10).Value = 23811404  # J136
$ws.Cells.Item(136, 11).Value = 5796.4998  # K136
$ws.Cells.Item(136, 12).Value = 71434212  # L136
$ws.Cells.Item(136, 13).Value = -3246.4998  # M136
$ws.Cells.Item(136, 14).Value = -71439312  # N136

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 3403.55  # H81
$ws.Cells.Item(81, 9).Value = 4067  # I81
$ws.Cells.Item(81, 10).Value = 2740.1  # J81
$ws.Cells.Item(81, 11).Value = 8134  # K81
$ws.Cells.Item(81, 12).Value = 5480.2  # L81
$ws.Cells.Item(81, 13).Value = -7073  # M81
$ws.Cells.Item(81, 14).Value = -7602.2  # N81
# Row 84
$ws.Cells.Item(84, 8).Value = 3403.55  # H84
$ws.Cells.Item(84, 9).Value = 4067  # I84
$ws.Cells.Item(84, 10).Value = 2740.1  # J84
$ws.Cells.Item(84, 11).Value = 40670  # K84
$ws.Cells.Item(84, 12).Value = 27401  # L84
$ws.Cells.Item(84, 13).Value = -35366  # M84
$ws.Cells.Item(84, 14).Value = -38009  # N84
# Row 107
$ws.Cells.Item(107, 8).Value = 1199.5714  # H107
$ws.Cells.Item(107, 9).Value = 1232.8334  # I107
$ws.Cells.Item(107, 10).Value = 1000  # J107
$ws.Cells.Item(107, 11).Value = 3698.5002  # K107
$ws.Cells.Item(107, 12).Value = 3000  # L107
$ws.Cells.Item(107, 13).Value = -1778.5002  # M107
$ws.Cells.Item(107, 14).Value = -6840  # N107
# Row 123
$ws.Cells.Item(123, 8).Value = 30085.8  # H123
$ws.Cells.Item(123, 10).Value = 45214.5  # J123
$ws.Cells.Item(123, 12).Value = 45214.5  # L123
$ws.Cells.Item(123, 14).Value = -55014.5  # N123
# Row 132
$ws.Cells.Item(132, 8).Value = 3207394.2  # H132
$ws.Cells.Item(132, 9).Value = 2763.6538  # I132
$ws.Cells.Item(132, 10).Value = 6412025  # J132
$ws.Cells.Item(132, 11).Value = 8290.9614  # K132
$ws.Cells.Item(132, 12).Value = 19236075  # L132
$ws.Cells.Item(132, 13).Value = -5760.9614  # M132
$ws.Cells.Item(132, 14).Value = -19241135  # N132
# Row 136
$ws.Cells.Item(136, 8).Value = 2043.3788  # H136
$ws.Cells.Item(136, 9).Value = 1705.72  # I136
$ws.Cells.Item(136, 10).Value = 3098.5625  # J136
$ws.Cells.Item(136, 11).Value = 5117.16  # K136
$ws.Cells.Item(136, 12).Value = 9295.6875  # L136
$ws.Cells.Item(136, 13).Value = -2567.16  # M136
$ws.Cells.Item(136, 14).Value = -14395.6875  # N136
